$d = $word.ActiveDocument

$old = "Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός Ωρίωνα 2022: 16-25 Ιανουαρίου, 14-23 Φεβρουαρίου, 14-24 Μαρτίου"
$new = "2022 Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός Ωρίωνα: 16-25 Ιανουαρίου, 14-23 Φεβρουαρίου, 14-24 Μαρτίου"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
